$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("tables")
$ws1.Rows.Item(7).Delete()
Write-Output $ws1.UsedRange.Address()
